$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# IBGE states table: existing single data row ("baiano"/Salvador/...) moves
# down to row 6, and four new rows are inserted above it (rows 2-5).
$data = @(
    @("acriano",    "Rio Branco", "GLADSON DE LIMA CAMELI",         "906.876",    "0,663"),
    @("alagoano",   "Maceió",     "PAULO SURUAGY DO AMARAL DANTAS", "3.365.351",  "0,631"),
    @("amapaense",  "Macapá",     "CLÉCIO LUÍS VILHENA VIEIRA",     "877.613",    "0,708"),
    @("amazonense", "Manaus",     "WILSON MIRANDA LIMA",            "4.269.995",  "0,674"),
    @("baiano",     "Salvador",   "JERÔNIMO RODRIGUES SOUZA",       "14.985.284", "0,660")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]

    # Columns D (Populacao) and E (IDH) hold digit/punctuation strings
    # ("906.876", "0,663", ...) that Excel's input parser would otherwise
    # coerce into numbers. Write them through a literal-string formula and
    # flatten the formula to its computed value via copy/paste-special, so
    # the cell ends up holding plain text with no formula and no special
    # number-format styling left behind.
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Formula = '="' + $entry[3] + '"'
    $dCell.Copy()
    $dCell.PasteSpecial(-4163)

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Formula = '="' + $entry[4] + '"'
    $eCell.Copy()
    $eCell.PasteSpecial(-4163)

    $row++
}

$excel.CutCopyMode = 0
